# Append the latest Argent price observation as a new row at the bottom
# of the price history table (row 89): date "2025-01-26", value "5.83".
#
# Both the date and the value in this sheet are stored as literal text
# (not real Excel dates / numbers), matching every other row below the
# header. Excel's COM `.Value` setter auto-detects date- and
# number-looking strings and would otherwise coerce them to a serial
# date / numeric value, so we briefly force Text number formatting
# before writing the values, then restore the cell style to the
# workbook's default ("Normal") so the new cells don't end up carrying
# any extra formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1  # 88 existing data/header rows -> 89

$targetRange = $ws.Range("A" + $newRow + ":B" + $newRow)
$targetRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-01-26"
$ws.Cells.Item($newRow, 2).Value = "5.83"

$targetRange.Style = "Normal"
